$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "14_011215_0718_7_xgboost_binary_logits_with_random_3in1_preprocess_valid1_valid2_"
$ws.Range("B15").Value = 0.63186
$ws.Range("C15").Value = "ensembled 7 tree xgboost binary logits on random combined 3in1 data set with features preprocessed, with 2 valid sets"

$ws.Range("A16").Value = "15_011215_0818_7_tree_xgboost_binary_logits_and_1_linear_xgboost_binary_logits_with_random_3in1_preprocess_valid1_valid2_"
$ws.Range("B16").Value = 0.63528
$ws.Range("C16").Value = "ensembled 7 tree xgboost binary logits and 1 linear xgboost binary logits on random combined 3in1 data set with features preprocessed, with 2 valid sets"

$ws.Range("G17").Select()
